{"js": "// The document contains a single table of two-digit \u00f7 one-digit division\n// problems (\"NN\u00f7N=\") laid out 5-per-row, with blank spacer rows between\n// each populated row. This edit swaps in a fresh batch of problems while\n// preserving every run's formatting (font, size, alignment), by writing\n// into each table cell's `.value` (which replaces the cell's Range.Text in\n// place) rather than touching runs/paragraphs directly.\n//\n// Replacements are applied strictly in reading order (row-major, left to\n// right) because some source expressions repeat (e.g. \"20\u00f72=\" appears\n// twice) but map to different replacements depending on position, so a\n// global text find/replace would be ambiguous.\n\nconst replacements = [\n  \"19\u00f79=\", \"21\u00f73=\", \"59\u00f73=\", \"53\u00f76=\", \"89\u00f72=\",\n  \"36\u00f72=\", \"12\u00f74=\", \"39\u00f79=\", \"15\u00f79=\", \"46\u00f75=\",\n  \"29\u00f73=\", \"80\u00f77=\", \"44\u00f78=\", \"93\u00f73=\", \"58\u00f73=\",\n  \"40\u00f73=\", \"66\u00f78=\", \"92\u00f74=\", \"70\u00f78=\", \"77\u00f74=\",\n  \"85\u00f75=\", \"91\u00f76=\", \"70\u00f78=\", \"56\u00f74=\", \"62\u00f77=\",\n];\n\nconst expected = [\n  \"20\u00f72=\", \"44\u00f73=\", \"10\u00f75=\", \"20\u00f72=\", \"71\u00f76=\",\n  \"26\u00f77=\", \"77\u00f77=\", \"82\u00f73=\", \"10\u00f76=\", \"79\u00f77=\",\n  \"71\u00f79=\", \"56\u00f72=\", \"91\u00f79=\", \"53\u00f76=\", \"34\u00f75=\",\n  \"64\u00f72=\", \"11\u00f75=\", \"57\u00f74=\", \"60\u00f79=\", \"86\u00f78=\",\n  \"22\u00f72=\", \"38\u00f76=\", \"61\u00f77=\", \"88\u00f79=\", \"13\u00f73=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet k = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const current = grid[r][c];\n    if (current === \"\") continue;\n    if (k >= replacements.length) continue;\n    // Sanity-check we're replacing the expected cell before overwriting it.\n    if (current === expected[k]) {\n      table.getCell(r, c).value = replacements[k];\n    }\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document body holds a single table of two-digit / one-digit division\n# problems (\"NN\u00f7N=\") laid out 5 cells per row, with blank spacer rows\n# between each populated row. This script swaps in a fresh batch of\n# problems while preserving every cell's run formatting (font, size,\n# alignment) by assigning straight to `Cell.Range.Text`, which replaces\n# just the text content and keeps the existing run/paragraph properties.\n#\n# Replacements are applied strictly in reading order (row-major, left to\n# right) because some source expressions repeat (e.g. \"20\u00f72=\" appears\n# twice in row 1) but map to different replacements depending on\n# position, so a global text find/replace would be ambiguous.\n\n$expected = @(\n  \"20\u00f72=\", \"44\u00f73=\", \"10\u00f75=\", \"20\u00f72=\", \"71\u00f76=\",\n  \"26\u00f77=\", \"77\u00f77=\", \"82\u00f73=\", \"10\u00f76=\", \"79\u00f77=\",\n  \"71\u00f79=\", \"56\u00f72=\", \"91\u00f79=\", \"53\u00f76=\", \"34\u00f75=\",\n  \"64\u00f72=\", \"11\u00f75=\", \"57\u00f74=\", \"60\u00f79=\", \"86\u00f78=\",\n  \"22\u00f72=\", \"38\u00f76=\", \"61\u00f77=\", \"88\u00f79=\", \"13\u00f73=\"\n)\n\n$replacements = @(\n  \"19\u00f79=\", \"21\u00f73=\", \"59\u00f73=\", \"53\u00f76=\", \"89\u00f72=\",\n  \"36\u00f72=\", \"12\u00f74=\", \"39\u00f79=\", \"15\u00f79=\", \"46\u00f75=\",\n  \"29\u00f73=\", \"80\u00f77=\", \"44\u00f78=\", \"93\u00f73=\", \"58\u00f73=\",\n  \"40\u00f73=\", \"66\u00f78=\", \"92\u00f74=\", \"70\u00f78=\", \"77\u00f74=\",\n  \"85\u00f75=\", \"91\u00f76=\", \"70\u00f78=\", \"56\u00f74=\", \"62\u00f77=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $raw = $cell.Range.Text\n    $clean = $raw.TrimEnd([char]13, [char]7)\n    if ($clean -eq \"\") {\n      continue\n    }\n    if ($k -lt $replacements.Length -and $clean -eq $expected[$k]) {\n      $cell.Range.Text = $replacements[$k]\n    }\n    $k++\n  }\n}\n"}
